# Boolean Exempt Process Emissions From Carbon Tax (BEPEfCT)
# The boolean control lever, and its paired policy lever, can now be set
# separately for each industry sector instead of once globally.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BEPEfCT")

# List of industry sectors (A column, rows 2-26); each gets its own
# boolean setting in column B (defaulting to 0, same as the prior single
# global value).
$industries = @(
  "agriculture and forestry 01T03",
  "coal mining 05",
  "oil and gas extraction 06",
  "other mining and quarrying 07T08",
  "food beverage and tobacco 10T12",
  "textiles apparel and leather 13T15",
  "wood products 16",
  "pulp paper and printing 17T18",
  "refined petroleum and coke 19",
  "chemicals 20",
  "rubber and plastic products 22",
  "glass and glass products 231",
  "cement and other nonmetallic minerals 239",
  "iron and steel 241",
  "other metals 242",
  "metal products except machinery and vehicles 25",
  "computers and electronics 26",
  "appliances and electrical equipment 27",
  "other machinery 28",
  "road vehicles 29",
  "nonroad vehicles 30",
  "other manufacturing 31T33",
  "energy pipelines and gas processing 352T353",
  "water and waste 36T39",
  "construction 41T43"
)

# Row 2 used to hold the single word "Boolean" in A2 -- replace it (and
# the rows below) with one row per industry, each defaulting to 0.
for ($i = 0; $i -lt $industries.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $industries[$i]
    $ws.Range("B$row").Value = 0
}

# Row 1: keep the header label in B1, and add a new italic "units" label
# in A1.
$ws.Range("A1").Value = "Unit: boolean (0 or 1)"
$ws.Range("A1").Font.Italic = $true

# Widen column A so the longer industry names are readable.
$ws.Columns.Item(1).ColumnWidth = 46.33

# Restore a portrait page setup (matches the About sheet) now that the
# sheet has grown past a single screen.
$ws.PageSetup.Orientation = 1
